$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all target cells so values are preserved exactly
# (avoids Excel auto-converting numeric-looking / date-looking strings)
$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "B7", "C7", "D7", "E7", "B8", "C8", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "B20", "C20", "D20", "E20", "B21", "C21", "D21", "E21", "D22", "E22", "D23", "E23", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "E34", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "E44", "E45", "D46", "E46", "B47", "C47", "D47", "E47", "B48", "C48", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($c in $cells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply new values
$ws.Range('D2').Value = '79.558.46'
$ws.Range('E2').Value = '  +4.51%  '
$ws.Range('D3').Value = '3.198.15'
$ws.Range('E3').Value = '  +5.65%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '209.48'
$ws.Range('E5').Value = '  +6.28%  '
$ws.Range('D6').Value = '637.39'
$ws.Range('E6').Value = '  +2.65%  '
$ws.Range('B7').Value = 'Dogecoin'
$ws.Range('C7').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D7').Value = '0.253'
$ws.Range('E7').Value = '  +23.55%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.600'
$ws.Range('E9').Value = '  +9.78%  '
$ws.Range('D10').Value = '3.194.08'
$ws.Range('E10').Value = '  +5.61%  '
$ws.Range('D11').Value = '0.611'
$ws.Range('E11').Value = '  +39.73%  '
$ws.Range('D12').Value = '0.0000260'
$ws.Range('E12').Value = '  +36.02%  '
$ws.Range('E13').Value = '  +3.38%  '
$ws.Range('D14').Value = '5.41'
$ws.Range('E14').Value = '  +3.59%  '
$ws.Range('D15').Value = '3.787.08'
$ws.Range('E15').Value = '  +5.70%  '
$ws.Range('D16').Value = '32.22'
$ws.Range('E16').Value = '  +11.95%  '
$ws.Range('D17').Value = '79.498.88'
$ws.Range('E17').Value = '  +4.46%  '
$ws.Range('D18').Value = '3.191.44'
$ws.Range('E18').Value = '  +5.47%  '
$ws.Range('D19').Value = '14.56'
$ws.Range('E19').Value = '  +8.28%  '
$ws.Range('B20').Value = 'SuiNetwork'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D20').Value = '3.02'
$ws.Range('E20').Value = '  +28.15%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '9.36'
$ws.Range('E21').Value = '  +5.10%  '
$ws.Range('D22').Value = '441.90'
$ws.Range('E22').Value = '  +16.67%  '
$ws.Range('D23').Value = '5.25'
$ws.Range('E23').Value = '  +20.70%  '
$ws.Range('E24').Value = '  +12.40%  '
$ws.Range('D25').Value = '3.369.53'
$ws.Range('E25').Value = '  +6.31%  '
$ws.Range('D26').Value = '77.24'
$ws.Range('E26').Value = '  +6.88%  '
$ws.Range('D27').Value = '10.88'
$ws.Range('E27').Value = '  +11.86%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.13%  '
$ws.Range('E29').Value = '  +15.51%  '
$ws.Range('D30').Value = '9.16'
$ws.Range('E30').Value = '  +11.35%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = '560.51'
$ws.Range('E32').Value = '  +14.29%  '
$ws.Range('D33').Value = '1.52'
$ws.Range('E33').Value = '  +9.89%  '
$ws.Range('E34').Value = '  +32.45%  '
$ws.Range('E35').Value = '  +6.43%  '
$ws.Range('D36').Value = '23.13'
$ws.Range('E36').Value = '  +12.87%  '
$ws.Range('D37').Value = '0.121'
$ws.Range('E37').Value = '  +17.86%  '
$ws.Range('D38').Value = '0.998'
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('D39').Value = '0.411'
$ws.Range('E39').Value = '  +8.98%  '
$ws.Range('D40').Value = '162.97'
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').Value = '20.25'
$ws.Range('E41').Value = '  +1.12%  '
$ws.Range('D42').Value = '5.65'
$ws.Range('E42').Value = '  +11.38%  '
$ws.Range('D43').Value = '192.70'
$ws.Range('E43').Value = '  +1.67%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E45').Value = '  +12.37%  '
$ws.Range('D46').Value = '0.800'
$ws.Range('E46').Value = '  +2.35%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '2.68'
$ws.Range('E47').Value = '  +10.90%  '
$ws.Range('B48').Value = 'ImmutableX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D48').Value = '1.34'
$ws.Range('E48').Value = '  +8.31%  '
$ws.Range('D49').Value = '42.82'
$ws.Range('E49').Value = '  +3.15%  '
$ws.Range('D50').Value = '25.98'
$ws.Range('E50').Value = '  +17.76%  '
$ws.Range('D51').Value = '0.639'
$ws.Range('E51').Value = '  +7.50%  '

# Restore default style (remove text-format marker / quote-prefix artifacts)
foreach ($c in $cells) {
    $ws.Range($c).Style = "Normal"
}
